# "Generate Report for Handoff"
#
# The localization-status report is refreshed: the Overview / per-locale
# "Status" cells flip from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated "Latest Handoff"/"Latest HO
# Xliff Generate Date" timestamps are bumped to reflect the new handoff
# run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (table columns: zh-cn, de-de, Latest HO Xliff Generate Date) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 09:01:35"

# --- zh-cn sheet (Status, Latest Handoff Datetime) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 09:01:30"

# --- de-de sheet (Status, Latest Handoff Datetime) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 09:01:35"

# --- Column widths: the Status column on every sheet re-fits to the
# shorter "Ready for handoff" text. ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
